# Updates cryptos list price (D) and volume-change (E) columns
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a plain number
# (losing e.g. trailing zeros) must be marked as Text first, so the literal
# string is preserved exactly like the other inline-string cells in the sheet.
$ws.Range("D2").Value = "61.405.09"
$ws.Range("D3").Value = "2.993.57"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.35"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.66"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "2.988.27"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.98"
$ws.Range("D15").Value = "3.478.31"
$ws.Range("D17").Value = "61.428.27"
$ws.Range("D18").Value = "2.990.51"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.89"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.94"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.14"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.98"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.57"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.53"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.58"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "452.17"
$ws.Range("D38").Value = "3.155.90"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.96"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.00"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.94"
$ws.Range("D50").Value = "0.0₃0496"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"

# Volume(1h) column - percentage strings already contain surrounding spaces
# so Excel keeps them as text without any extra formatting tricks.
$ws.Range("E2").Value = "  -3.68%  "
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -5.00%  "
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("E20").Value = "  -4.62%  "
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("E28").Value = "  -6.40%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  -8.57%  "
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -7.75%  "
$ws.Range("E44").Value = "  +8.29%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -5.27%  "
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("E50").Value = "  -8.00%  "
$ws.Range("E51").Value = "  +6.12%  "

Write-Host "Updated cryptos list prices and volume changes"
